$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Data changes ---
$ws1.Range("A2").Value = "LoginPage_TC1"
$ws1.Range("B2").Value = "loginToApp"

$ws2.Range("A2").Value = "LoginPage_TC1"
$ws2.Range("B2").Value = "'tejeshkumar.gangari@gmail.com"
$ws2.Range("C2").Value = "'Test@123"
$ws2.Range("B2:C2").HorizontalAlignment = -4108

# --- View / selection changes ---
$ws2.Range("F8").Select()
$ws1.Activate()
$ws1.Range("D7").Select()
